# Insert a new data row at row 136 (pushing the former rows 136-210 down to 137-211)
# and populate it with the new price-record values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(136).Insert()

$ws.Range("A136").Value2 = 11
$ws.Range("B136").Value2 = "Vega Monumental Concepción"
$ws.Range("C136").Value2 = "Bíobío"
$ws.Range("D136").Value2 = 44582
$ws.Range("E136").Value2 = 8
$ws.Range("F136").Value2 = "Fruta"
$ws.Range("G136").Value2 = 100102
$ws.Range("H136").Value2 = "Cítricos"
$ws.Range("I136").Value2 = 100102005
$ws.Range("J136").Value2 = "Naranja"
$ws.Range("K136").Value2 = "Valencia"
$ws.Range("L136").Value2 = "Primera"
$ws.Range("M136").Value2 = 290
$ws.Range("N136").Value2 = 9000
$ws.Range("O136").Value2 = 9500
$ws.Range("P136").Value2 = 9241
$ws.Range("Q136").Value2 = "$/bandeja 15 kilos granel"
$ws.Range("R136").Value2 = "Región de O'Higgins"
$ws.Range("S136").Value2 = 616
$ws.Range("T136").Value2 = 15
